# "Generate Report for Handback"
#
# The handback-status report records, for each localized file, the
# timestamps of the handoff/handback xliff generation steps. This run
# regenerated the report, advancing a handful of the recorded
# timestamps (the rest of the workbook content is unchanged):
#
#   Overview!G2            (Latest HO Xliff Generate Date)      07:02:29 -> 07:03:23
#   de-de!H2                (Correspond Handoff Datetime)        07:02:29 -> 07:03:23  (same run as above)
#   zh-cn!H2                (Correspond Handoff Datetime)        07:02:24 -> 07:03:19
#   zh-cn!K2                (Correspond Handback DateTime)       07:02:48 -> 07:03:39
#   de-de!K2                (Correspond Handback DateTime)       07:02:55 -> 07:03:45

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first (9a6ee2a4...) entry
$wsOverview.Range("G2").Value = "2016-08-22 07:03:23"

# zh-cn sheet: Correspond Handoff / Handback datetimes for the first entry
$wsZhCn.Range("H2").Value = "2016-08-22 07:03:19"
$wsZhCn.Range("K2").Value = "2016-08-22 07:03:39"

# de-de sheet: Correspond Handoff / Handback datetimes for the first entry
$wsDeDe.Range("H2").Value = "2016-08-22 07:03:23"
$wsDeDe.Range("K2").Value = "2016-08-22 07:03:45"
